# "Generate Report for Handback" - update the localization status report to
# reflect that the handback for file a.md has completed (de-de and zh-cn).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Overview")
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws3 = $wb.Worksheets.Item("de-de")

$aMdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9bc174efeafc1b3cbe9bf2d2d5e5af01467ad0d2/e2e/a.md"
$bMdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9bc174efeafc1b3cbe9bf2d2d5e5af01467ad0d2/e2e/b.md"

# ---------------------------------------------------------------------------
# Status text: "Ready for handoff" -> "Handed back: in sync with en-US"
# (Overview summary columns + the Status column on each language sheet)
# ---------------------------------------------------------------------------
$handedBack = "Handed back: in sync with en-US"

$ws1.Range("E2").Value = $handedBack
$ws1.Range("F2").Value = $handedBack
$ws1.Range("E3").Value = $handedBack
$ws1.Range("F3").Value = $handedBack

$ws2.Range("C2").Value = $handedBack
$ws2.Range("C3").Value = $handedBack

$ws3.Range("C2").Value = $handedBack
$ws3.Range("C3").Value = $handedBack

# ---------------------------------------------------------------------------
# zh-cn sheet (table1 / sheet2): Latest Target File (I) & Latest Handback
# File (J) now populated; Latest Handback DateTime (K) stamped.
# ---------------------------------------------------------------------------
$ws2.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), $aMdUrl, [Type]::Missing, [Type]::Missing, "a.md")
$ws2.Hyperlinks.Add($ws2.Range("I2"), $aMdUrl, [Type]::Missing, [Type]::Missing, "a.md")
$ws2.Hyperlinks.Add($ws2.Range("A3"), $bMdUrl, [Type]::Missing, [Type]::Missing, "b.md")
$ws2.Hyperlinks.Add($ws2.Range("I3"), $aMdUrl, [Type]::Missing, [Type]::Missing, "a.md")

$ws2.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$ws2.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$ws2.Range("K2").Value = "2016-09-07 08:49:00"
$ws2.Range("K3").Value = "2016-09-07 08:49:00"

# ---------------------------------------------------------------------------
# de-de sheet (table2 / sheet3): same pattern, different handback file /
# handback datetime.
# ---------------------------------------------------------------------------
$ws3.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), $aMdUrl, [Type]::Missing, [Type]::Missing, "a.md")
$ws3.Hyperlinks.Add($ws3.Range("I2"), $aMdUrl, [Type]::Missing, [Type]::Missing, "a.md")
$ws3.Hyperlinks.Add($ws3.Range("A3"), $bMdUrl, [Type]::Missing, [Type]::Missing, "b.md")
$ws3.Hyperlinks.Add($ws3.Range("I3"), $aMdUrl, [Type]::Missing, [Type]::Missing, "a.md")

$ws3.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$ws3.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$ws3.Range("K2").Value = "2016-09-07 08:49:27"
$ws3.Range("K3").Value = "2016-09-07 08:49:27"

# ---------------------------------------------------------------------------
# Column width adjustments (widened to fit the longer status / file-name
# text that is now displayed). The saved OOXML "width" attribute is derived
# from ColumnWidth with pixel-level rounding, so the inputs below are chosen
# to land on the desired stored widths (~29.98 and 40 characters).
# ---------------------------------------------------------------------------
$ws1.Range("E1").ColumnWidth = 29.15
$ws1.Range("F1").ColumnWidth = 29.15

$ws2.Range("C1").ColumnWidth = 29.15
$ws2.Range("J1").ColumnWidth = 39.15

$ws3.Range("C1").ColumnWidth = 29.15
$ws3.Range("J1").ColumnWidth = 39.15
